$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.118.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.910.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  -0.88%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4815'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3821'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07367'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9336'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07796'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.910.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("E14").Value = '  +1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.652'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008828'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.152.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.171'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.143.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.967'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08956'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.340'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.257'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7739'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.32%  '
$ws.Range("E35").Value = '  +1.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02060'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.109'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05329'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5500'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.990'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.028'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1530'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.486'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4842'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.14%  '
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.657'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06084'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
